$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc7 = "Unet Efficientnet-b1 `nTrained on full train data`n(padif needed + 224x224 resize)"
$log7 = "Epoch 1 | Train Loss: 0.213 | Valid Loss: 0.126 | Combined metric: 0.664 | Dice: 0.298 (LB 0.649, SB 0.492, S 0.633) | Hausdorff: 0.092 (LB 0.103, SB 0.091, S 0.081)`nEpoch 2 | Train Loss: 0.113 | Valid Loss: 0.136 | Combined metric: 0.659 | Dice: 0.272 (LB 0.663, SB 0.495, S 0.594) | Hausdorff: 0.083 (LB 0.096, SB 0.086, S 0.066)`nEpoch 3 | Train Loss: 0.098 | Valid Loss: 0.124 | Combined metric: 0.648 | Dice: 0.286 (LB 0.668, SB 0.513, S 0.651) | Hausdorff: 0.111 (LB 0.123, SB 0.111, S 0.099)`nEpoch 4 | Train Loss: 0.089 | Valid Loss: 0.140 | Combined metric: 0.650 | Dice: 0.269 (LB 0.619, SB 0.482, S 0.592) | Hausdorff: 0.095 (LB 0.103, SB 0.075, S 0.108)`nEpoch 5 | Train Loss: 0.086 | Valid Loss: 0.116 | Combined metric: 0.685 | Dice: 0.302 (LB 0.662, SB 0.563, S 0.706) | Hausdorff: 0.059 (LB 0.074, SB 0.063, S 0.040)"

$ws.Range("A7").Value = $desc7
$ws.Range("B7").Value = $log7
$ws.Range("C7").Value = 0.685
$ws.Range("D7").Value = 0.79655
$ws.Range("E7").Value = 0.78819

$ws.Rows.Item(7).RowHeight = 57.45

$ws.Range("A7:B7").WrapText = $true

$ws.Range("C10").Select()
